$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$updates_ALC = @{
    "H64" = 4162.6665
    "I64" = 4864.6665
    "J64" = 3811.6667
    "K64" = 4864.6665
    "L64" = 3811.6667
    "M64" = -4616.6665
    "N64" = -4307.6667
    "H67" = 4162.6665
    "I67" = 4864.6665
    "J67" = 3811.6667
    "K67" = 4864.6665
    "L67" = 3811.6667
    "M67" = -4006.6665
    "N67" = -5527.6667
    "H69" = 4000
    "J69" = 4000
    "L69" = 12000
    "N69" = -13748
    "H70" = 1503
    "I70" = 0
    "J70" = 1503
    "K70" = 0
    "L70" = 4509
    "N70" = -5049
    "H72" = 4000
    "J72" = 4000
    "L72" = 36000
    "N72" = -44736
    "H73" = 1503
    "I73" = 0
    "J73" = 1503
    "K73" = 0
    "L73" = 4509
    "N73" = -6381
    "H76" = 6345.636
    "J76" = 7166.8887
    "L76" = 7166.8887
    "N76" = -7796.8887
    "H79" = 6345.636
    "J79" = 7166.8887
    "L79" = 7166.8887
    "N79" = -9350.8887
    "H80" = 1264.4
    "I80" = 2800
    "J80" = 880.5
    "K80" = 8400
    "L80" = 2641.5
    "M80" = -7402
    "N80" = -4637.5
    "H83" = 1264.4
    "I83" = 2800
    "J83" = 880.5
    "K83" = 25200
    "L83" = 7924.5
    "M83" = -20208
    "N83" = -17908.5
    "H86" = 6007.4614
    "I86" = 6999.6665
    "J86" = 5709.8
    "K86" = 6999.6665
    "L86" = 5709.8
    "M86" = -5876.6665
    "N86" = -7955.8
    "H89" = 6007.4614
    "I89" = 6999.6665
    "J89" = 5709.8
    "K89" = 34998.3325
    "L89" = 28549
    "M89" = -29382.3325
    "N89" = -39781
    "H92" = 574.95
    "I92" = 447.3158
    "J92" = 3000
    "K92" = 447.3158
    "L92" = 3000
    "M92" = 800.6841999999999
    "N92" = -5496
    "H97" = 2450
    "J97" = 2450
    "L97" = 7350
    "N97" = -8342
    "H98" = 4317.706
    "I98" = 4317.706
    "K98" = 4317.706
    "M98" = -2819.706
    "H101" = 485.5
    "I101" = 425.6
    "J101" = 785
    "K101" = 1276.8
    "L101" = 2355
    "M101" = 345.1999999999998
    "N101" = -5599
    "H103" = 1696.5
    "I103" = 300
    "J103" = 1975.8
    "K103" = 900
    "L103" = 5927.4
    "M103" = -314
    "N103" = -7099.4
    "H106" = 8968.6875
    "I106" = 9670.929
    "K106" = 9670.929
    "M106" = -9039.929
    "H122" = 4317.706
    "I122" = 4317.706
    "K122" = 12953.118
    "M122" = -10503.118
}
foreach ($key in $updates_ALC.Keys) {
    Set-Cell $ws $key $updates_ALC[$key]
}
$ws.Range("M70").ClearContents()
$ws.Range("M73").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$updates_ARM = @{
    "H88" = 2771.3572
    "J88" = 2800
    "L88" = 2800
    "N88" = -3612
    "H91" = 2771.3572
    "J91" = 2800
    "L91" = 2800
    "N91" = -5608
    "H110" = 1191.762
    "I110" = 742.0714
    "K110" = 742.0714
    "M110" = 1302.9286
    "H132" = 2192.818
    "I132" = 2111.75
    "J132" = 2290.1
    "K132" = 6335.25
    "L132" = 6870.299999999999
    "M132" = -3805.25
    "N132" = -11930.3
}
foreach ($key in $updates_ARM.Keys) {
    Set-Cell $ws $key $updates_ARM[$key]
}

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$updates_CRP = @{
    "H62" = 40002100
    "I62" = 2500
    "J62" = 66668500
    "K62" = 2500
    "L62" = 66668500
    "M62" = -1876
    "N62" = -66669748
    "H65" = 40002100
    "I65" = 2500
    "J65" = 66668500
    "K65" = 12500
    "L65" = 333342500
    "M65" = -9380
    "N65" = -333348740
}
foreach ($key in $updates_CRP.Keys) {
    Set-Cell $ws $key $updates_CRP[$key]
}

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$updates_CUL = @{
    "H109" = 68870.734
    "I109" = 112140.11
    "J109" = 3966.6667
    "K109" = 336420.33
    "L109" = 11900.0001
    "M109" = -335380.33
    "N109" = -13980.0001
}
foreach ($key in $updates_CUL.Keys) {
    Set-Cell $ws $key $updates_CUL[$key]
}

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$updates_LTW = @{
    "H40" = 3079.8
    "I40" = 2867.6667
    "J40" = 3170.7144
    "K40" = 2867.6667
    "L40" = 3170.7144
    "M40" = -2731.6667
    "N40" = -3442.7144
    "H93" = 629.5
    "I93" = 629.5
    "K93" = 629.5
    "M93" = 618.5
    "H122" = 17865930
    "I122" = 25011550
    "J122" = 1877.5
    "K122" = 75034650
    "L122" = 5632.5
    "M122" = -75032200
    "N122" = -10532.5
}
foreach ($key in $updates_LTW.Keys) {
    Set-Cell $ws $key $updates_LTW[$key]
}
